$d = $word.ActiveDocument

# --- Paragraph 1: "Yerba Buena, 15 de Mayo de 1990" ---
$p1 = $d.Paragraphs.Item(1)
$p1.KeepWithNext = $true
$p1.SpaceAfter = 12

# --- Paragraph 2: "ORDENANZA Nº 367" -> bold, keepNext, spacing before/after ---
$p2 = $d.Paragraphs.Item(2)
$p2.KeepWithNext = $true
$p2.SpaceBefore = 12
$p2.SpaceAfter = 18
$p2.Range.Font.Bold = 1

# --- Paragraph 3: "EL " + "CONCEJO..." merge into one run, bold, centered, indented ---
$rng3 = $d.Content
$rng3.Find.Execute("EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA", $false, $false, $false, $false, $false, $true, 1, $false, "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA", 2) | Out-Null

$p3 = $d.Paragraphs.Item(3)
$p3.KeepWithNext = $true
$p3.SpaceBefore = 18
$p3.SpaceAfter = 18
$p3.LeftIndent = 99.2
$p3.RightIndent = 99.2
$p3.Range.Font.Bold = 1

# --- Paragraph 4: "ARTICULO PRIMERO: " -> keepNext, spacing after, underline split ---
$p4 = $d.Paragraphs.Item(4)
$p4.KeepWithNext = $true
$p4.SpaceAfter = 6

$rngA = $d.Content
$rngA.Find.Execute("ARTICULO PRIMERO", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngA.Font.Underline = 1

$rngB = $d.Content
$rngB.Find.Execute("ARTICULO PRIMERO:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$colonB = $d.Range($rngB.End - 1, $rngB.End)
$colonB.Font.Underline = 1

# --- Paragraph 5: "ARTICULO SEGUNDO: " -> same treatment ---
$p5 = $d.Paragraphs.Item(5)
$p5.KeepWithNext = $true
$p5.SpaceAfter = 6

$rngC = $d.Content
$rngC.Find.Execute("ARTICULO SEGUNDO", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngC.Font.Underline = 1

$rngD = $d.Content
$rngD.Find.Execute("ARTICULO SEGUNDO:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$colonD = $d.Range($rngD.End - 1, $rngD.End)
$colonD.Font.Underline = 1

# --- Paragraph 6: "ARTICULO TERCERO: " -> same treatment ---
$p6 = $d.Paragraphs.Item(6)
$p6.KeepWithNext = $true
$p6.SpaceAfter = 6

$rngE = $d.Content
$rngE.Find.Execute("ARTICULO TERCERO", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngE.Font.Underline = 1

$rngF = $d.Content
$rngF.Find.Execute("ARTICULO TERCERO:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$colonF = $d.Range($rngF.End - 1, $rngF.End)
$colonF.Font.Underline = 1

# --- Section: footer + starting page number ---
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftr.Range.ParagraphFormat.Style = "Footer"
$ftr.Range.Font.NameAscii = "Book Antiqua"
$ftr.Range.Font.Name = "Book Antiqua"
$ftr.Range.Font.Size = 10
$ftr.Range.Font.TextColor.ObjectThemeColor = 12
$ftr.PageNumbers.StartingNumber = 260

Write-Output "done"
